$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new values look like plain numbers to Excel need to be
# forced to Text format first, so they are stored as text (matching the source
# data which is all text), not auto-converted to numeric values.
$textCells = @("D4", "D5", "D6", "D9", "D11", "D15", "D20", "D22", "D25", "D28", "D30", "D31", "D32", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D45", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the diff.
$ws.Range("D2").Value = '65.225.71'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '3.568.57'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '600.69'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '133.70'
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").Value = '3.566.34'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").Value = '4.157.13'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("E14").Value = '  -2.59%  '
$ws.Range("D15").Value = '27.06'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").Value = '3.555.57'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '64.478.72'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").Value = '14.48'
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '391.92'
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("D24").Value = '3.706.12'
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '74.26'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").Value = '7.85'
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").Value = '  +25.33%  '
$ws.Range("D30").Value = '8.63'
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '2.30'
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").Value = '3.566.15'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '0.147'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Value = '170.71'
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("D38").Value = '6.97'
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '5.08'
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = '1.54'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = '0.0816'
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = '26.60'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  +5.39%  '
$ws.Range("D45").Value = '43.08'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").Value = '6.93'
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").Value = '2.445.42'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '0.0266'
$ws.Range("E51").Value = '  +0.28%  '
